$p = $ppt.ActivePresentation

# Slide 10 ("7. Child component updating parent's state" -> new title about
# parent <-> child write-back data flow)
$slide10 = $p.Slides.Item(10)
$title10 = $slide10.Shapes.Item(1)
$title10.TextFrame.TextRange.Text = "7. Parent component <-> Child. Data, “write-back”"

# Slide 15 ("11. Routing" -> "11. Routing basics")
$slide15 = $p.Slides.Item(15)
$title15 = $slide15.Shapes.Item(1)
$title15.TextFrame.TextRange.Text = "11. Routing basics"
